$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1099.7407
$ws.Range("I28").Value = 376.65
$ws.Range("K28").Value = 376.65
$ws.Range("M28").Value = 108.35
$ws.Range("H34").Value = 9451.5
$ws.Range("I34").Value = 9451.5
$ws.Range("K34").Value = 9451.5
$ws.Range("M34").Value = -9248.5
$ws.Range("H36").Value = 9451.5
$ws.Range("I36").Value = 9451.5
$ws.Range("K36").Value = 9451.5
$ws.Range("M36").Value = -8736.5
$ws.Range("H40").Value = 4314.2593
$ws.Range("I40").Value = 2483.4285
$ws.Range("J40").Value = 4955.05
$ws.Range("K40").Value = 2483.4285
$ws.Range("L40").Value = 4955.05
$ws.Range("M40").Value = -2308.4285
$ws.Range("N40").Value = -5305.05
$ws.Range("H97").Value = 2050
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H112").Value = 4653.3228
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H132").Value = 29415408
$ws.Range("I132").Value = 33337064
$ws.Range("K132").Value = 100011192
$ws.Range("M132").Value = -100008662
$ws.Range("H137").Value = 52540.086
$ws.Range("I137").Value = 82039.23
$ws.Range("J137").Value = 2618.4614
$ws.Range("K137").Value = 246117.69
$ws.Range("L137").Value = 7855.3842
$ws.Range("M137").Value = -243567.69
$ws.Range("N137").Value = -12955.3842
$ws.Range("H138").Value = 3951.8372
$ws.Range("I138").Value = 6000
$ws.Range("J138").Value = 3851.9268
$ws.Range("K138").Value = 18000
$ws.Range("L138").Value = 11555.7804
$ws.Range("M138").Value = -12860
$ws.Range("N138").Value = -21835.7804

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 22398.531
$ws.Range("I74").Value = 3267.6858
$ws.Range("K74").Value = 3267.6858
$ws.Range("M74").Value = -2393.6858
$ws.Range("H77").Value = 22398.531
$ws.Range("I77").Value = 3267.6858
$ws.Range("K77").Value = 16338.429
$ws.Range("M77").Value = -11970.429
$ws.Range("H110").Value = 3473786.8
$ws.Range("I110").Value = 3969613.5
$ws.Range("K110").Value = 3969613.5
$ws.Range("M110").Value = -3967568.5
$ws.Range("H132").Value = 2556.4
$ws.Range("I132").Value = 1667.92
$ws.Range("K132").Value = 5003.76
$ws.Range("M132").Value = -2473.76

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3584453.2
$ws.Range("I94").Value = 4550004
$ws.Range("K94").Value = 4550004
$ws.Range("M94").Value = -4549553
$ws.Range("H107").Value = 3973412
$ws.Range("I107").Value = 5107709.5
$ws.Range("J107").Value = 3370.75
$ws.Range("K107").Value = 5107709.5
$ws.Range("L107").Value = 3370.75
$ws.Range("M107").Value = -5105789.5
$ws.Range("N107").Value = -7210.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 131.9
$ws.Range("I7").Value = 131.9
$ws.Range("K7").Value = 131.9
$ws.Range("M7").Value = -18.90000000000001
$ws.Range("H58").Value = 7284.84
$ws.Range("I58").Value = 8634.875
$ws.Range("K58").Value = 8634.875
$ws.Range("M58").Value = -8431.875
$ws.Range("H88").Value = 72697
$ws.Range("J88").Value = 72697
$ws.Range("L88").Value = 72697
$ws.Range("N88").Value = -73509
$ws.Range("H91").Value = 72697
$ws.Range("J91").Value = 72697
$ws.Range("L91").Value = 72697
$ws.Range("N91").Value = -75505
$ws.Range("H103").Value = 9680.200000000001
$ws.Range("I103").Value = 2441.3333
$ws.Range("J103").Value = 20538.5
$ws.Range("K103").Value = 2441.3333
$ws.Range("L103").Value = 20538.5
$ws.Range("M103").Value = -1269.3333
$ws.Range("N103").Value = -22882.5
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954
$ws.Range("H132").Value = 113736.234
$ws.Range("I132").Value = 92698
$ws.Range("J132").Value = 152306.33
$ws.Range("K132").Value = 278094
$ws.Range("L132").Value = 456918.99
$ws.Range("M132").Value = -275564
$ws.Range("N132").Value = -461978.99
$ws.Range("H136").Value = 7284.84
$ws.Range("I136").Value = 8634.875
$ws.Range("K136").Value = 25904.625
$ws.Range("M136").Value = -23354.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 46254.453
$ws.Range("J37").Value = 46254.453
$ws.Range("L37").Value = 138763.359
$ws.Range("N37").Value = -138987.359
$ws.Range("H64").Value = 1300
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1300
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3900
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4440
$ws.Range("H67").Value = 1300
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1300
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3900
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5772
$ws.Range("H98").Value = 1330.091
$ws.Range("J98").Value = 1373.1
$ws.Range("L98").Value = 4119.299999999999
$ws.Range("N98").Value = -7115.299999999999
$ws.Range("H109").Value = 1189.7778
$ws.Range("I109").Value = 1189.7778
$ws.Range("K109").Value = 3569.3334
$ws.Range("M109").Value = -2529.3334
$ws.Range("H112").Value = 500
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H122").Value = 1410.4117
$ws.Range("J122").Value = 1633.5714
$ws.Range("L122").Value = 14702.1426
$ws.Range("N122").Value = -19602.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6136699
$ws.Range("I126").Value = 5684706
$ws.Range("K126").Value = 17054118
$ws.Range("M126").Value = -17051648
$ws.Range("H132").Value = 3277
$ws.Range("I132").Value = 3202.3044
$ws.Range("K132").Value = 9606.913199999999
$ws.Range("M132").Value = -7076.913199999999
$ws.Range("H136").Value = 13097.954
$ws.Range("J136").Value = 13097.954
$ws.Range("L136").Value = 39293.862
$ws.Range("N136").Value = -44393.862

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 34237.926
$ws.Range("I22").Value = 56458.5
$ws.Range("J22").Value = 1917.091
$ws.Range("K22").Value = 56458.5
$ws.Range("L22").Value = 1917.091
$ws.Range("M22").Value = -56163.5
$ws.Range("N22").Value = -2507.091
$ws.Range("H27").Value = 34237.926
$ws.Range("I27").Value = 56458.5
$ws.Range("J27").Value = 1917.091
$ws.Range("K27").Value = 56458.5
$ws.Range("L27").Value = 1917.091
$ws.Range("M27").Value = -56351.5
$ws.Range("N27").Value = -2131.091
$ws.Range("H42").Value = 11508.667
$ws.Range("J42").Value = 11508.667
$ws.Range("L42").Value = 11508.667
$ws.Range("N42").Value = -12634.667
$ws.Range("H46").Value = 1816559.5
$ws.Range("I46").Value = 43478260
$ws.Range("K46").Value = 43478260
$ws.Range("M46").Value = -43478072
$ws.Range("H49").Value = 11508.667
$ws.Range("J49").Value = 11508.667
$ws.Range("L49").Value = 11508.667
$ws.Range("N49").Value = -11802.667
$ws.Range("H100").Value = 2568.5715
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 2580
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 2580
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -3662

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1550.3462
$ws.Range("J113").Value = 2938.6
$ws.Range("L113").Value = 8815.799999999999
$ws.Range("N113").Value = -13155.8
$ws.Range("H132").Value = 17432350
$ws.Range("I132").Value = 23259634
$ws.Range("K132").Value = 69778902
$ws.Range("M132").Value = -69776372
